# Apply the "Added User Story descrip docx" edit.
#
# The change inserts one new row above the existing "User Guide / Manual"
# row (sheet row 81) of the "MyBar Artifacts" table, shifting every row
# below it down by one. The new row documents a "User Stories Description"
# artifact contributed by "Adam Clark". The underlying Excel Table
# (Table1) is grown by one row and the frozen-pane selection state is
# updated to reflect where the user was last working in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new blank row before row 81 - this shifts rows 81:172 down
#    to 82:173 and keeps their values/styles intact.
$ws.Rows.Item(81).Insert()

# 2. Grow Table1 (the ListObject covering A2:F111) so it covers the newly
#    inserted row as well (A2:F112).
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A2:F112"))

# 3. Fill in the new row with the "User Stories Description" artifact.
$ws.Range("A81").Value = ">"
$ws.Range("B81").Value = "User Stories Description"
$ws.Range("C81").Value = "Adam Clark"

# 4. Update the view state: the sheet was scrolled/selected differently
#    after the edit (frozen-pane scroll position and active selection).
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 3
$ws.Range("C82").Select()
